$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching the style of the other header cells (B1:G1)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data column values for "Save"
$saveValues = @(0, 0, 1, 0, 0, 1, 0, 0)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
